$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14: 2025-09-09 21:20:57 EUR->ARS rate update.
# A14 ("2025-09-09") looks like an ISO date, so a plain .Value assignment
# would get auto-parsed into a date serial number. Force it to be stored
# as literal text (matching the original sheet's plain-text date column)
# by entering it with a leading apostrophe, then resetting the cell style
# back to Normal so no stray number-format/quote-prefix styling lingers.
$ws.Range("A14").Value = "'2025-09-09"
$ws.Range("A14").Style = "Normal"

$ws.Range("B14").Value = "21:20:57"
$ws.Range("C14").Value = "1.00 EUR = 1685.7993 ARS"
